$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 44 (Excel row 45) goes live - fill in Ser, Date, Ayats, Content, Author, Tags
$ws.Cells.Item(45, 1).Value = 44
$ws.Cells.Item(45, 2).Value = (Get-Date -Year 2020 -Month 2 -Day 13 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(45, 2).NumberFormat = "d-mmm-yy"
$ws.Cells.Item(45, 3).Value = 'Surah Al Nisa, 17 - 29'
$ws.Cells.Item(45, 6).Value = 'Normal vs Abnormal, Marrying Prisoner of Wars, Energy around human, Colors, Animal''s vision'
$ws.Cells.Item(45, 4).Value = 'h1: Maintaining sanity
p: Quran is a diverse book. It is difficult to put its verses in boxes. Verses are interlinked with each other, their revelation circumstances, present happenings and reader’s intent. Each aspect varies the meaning. It is not possible to describe the meaning of a verse in its totality from a single perspective. 
p: Understanding Quran is not just reading it, but acting upon it. Lack of action, makes us go wary off its meaning that leads us to stop believing in it. <b class=lavendar>It speaks of a colorful life where people can marry the captured prisoners, even if they are married. </b>
p: It looks abnormal and cruel. It is totally abnormal to take somebody else’s wife into marriage for personal needs. 
quote: And [also prohibited to you are all] married women except those your right hands possess. [This is] the decree of Allah upon you. And lawful to you are [all others] beyond these, [provided] that you seek them [in marriage] with [gifts from] your property, desiring chastity, not unlawful sexual intercourse. So for whatever you enjoy [of marriage] from them, give them their due compensation as an obligation. And there is no blame upon you for what you mutually agree to beyond the obligation. Indeed, Allah is ever Knowing and Wise. <br> - Surah Al Nisa verse 24
p: May be it is abnormal. But to understand the abnormality of these orders we need to find what is normal. 
h3: What is normal
p: Right now, I am sitting on the prayer rug beside a heater. My fingers tapping on this keyboard. Screen is black. Laptop is connected to the power supply. I am looking at the screen and thinking;
p.b-left: Heater is warm. Warm is a feeling. Warm is made up of a mesh of electric connections that is creating fire. Some neurons in my body are feeling the heat. How? 
p.b-left: Fingers are tapping on the keyboard. Sequence of keyboard keys are registered in my mind. Fingers are pressing down, precisely and completely. The pressing of the button to understanding millions of combinations that create useful sentences. How can I encompass this information at one moment of time. Each moment, trillionth part of a second, the possible combinations the keyboard can be triggered, changes. How does I do it?
p.b-left: The screen is black. Black is a color. When a dog sees black, it finds no difference from any-other color. Dogs and cats are color blind. No animal on this planet can see all colors as humans do. The only creature triggering keyboard to make useful meanings come out of laptop, is the only creature that can see colors. Trillion of combinations in each second are linked to humans. How do I name the black, the real black? 
p.b-left: Laptop and its connection to electricity is giving its battery the charge which is needed to keep it alive. Batteries convert chemical energy into electrical energy. The word<b> energy</b> is undefined and has been justified an infinite number of times. It can not be created or destroyed, it stays there. 
h3: So what is normal around us?
p: Nothing is normal. Each passing moment, we experience a magic and we reaffirm ourselves, it is normal. There is nothing normal in this world. It is all abnormal. The movement of earth, to sun throwing light at us, to being the only planet in the world, created from the minutest particles of neurones, protons and electrons, is all but not normal. 
h3: What if..
p.b-left: We are Adam. 
p.b-left: We are kept inside a closed room for the first 30 years of life and we do not see the world around.
p: We will not believe the normal world. The sun, flow of river, the bird, clouds with water, blood in veins, insects making lines, fishes making patterns. We will not believe any of these magics if we come to Earth from another planet or space. 
p: We have become habitual of living a life in a magical world. The moment the dimensions of our life changes, we get afraid of losing the part with which we are habitually connected. We are forgetting there is only 1 x Guy who is claiming the creation of all these wonders. 
p: We are <b>afraid-muslims</b> or <b>afraid-believers</b>, because we have become habitual to the magic of the world. Its dimensions needs to be altered every day. <b>Each day, we need to jump atleast once to find Him</b>. 
h3: Conclusion
p: Becoming a believer is only possible after we find the abnormal part of our life. We can not truly trust the God, until we keep ourselves on our feet and we keep our eyes open. The lesser we think, the more we go away from Him. 
p: If all the things happening around us all abnormal, there is no point of highlighting one more abnormal in Quran and labelling it weird. We may marry the captured girls from war, because it is said by Quran. The only book of 4 holy books that claims the connection to the creator of all the magics in this world.
p.note: Guys, if you like this project. Please follow this project''s page on twitter. <a href="https://twitter.com/zakatlists">Click here to go to the twitter page</a>. Subscribe below to receive updates.'
$ws.Cells.Item(45, 5).Value = "Qasim Ali"

$ws.Rows.Item(45).RowHeight = 409.6

$ws.Application.ActiveWindow.TopLeftCell = $ws.Range("A45")
$ws.Range("D45").Select()
